# Generate Report for Handoff
# Adds two new entries (21d5ef20-... and 4999b8bf-...) to the localization
# status report, both landing with status "Ready for handoff", just above
# the pre-existing "ac98edee-..." row (which itself stays "Ready for
# handoff" and keeps its own dates/links).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Duplicate the last data row (row 5, "ac98edee...") twice, pushing it to
# row 7 and opening up rows 5-6 (with matching style) for the new entries.
$ws1.Rows.Item(5).Copy()
$ws1.Rows.Item(6).Insert()
$ws1.Rows.Item(5).Copy()
$ws1.Rows.Item(6).Insert()

$ws1.Range("A5").Value = "21d5ef20-4a71-4957-bc19-8864fef6ab6f.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-23 22:40:20"

$ws1.Range("A6").Value = "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"
$ws1.Range("D6").Value = "2016-03-23 22:40:20"

# Row 7 retains the original ac98edee-... content that was duplicated down
# (no value changes needed there).

# Rebuild every hyperlink in this sheet, in final top-to-bottom order, so
# relationship ids come out sequential (rId2..rId7).
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a896c4a09a321d440863c6933282f9952176d25/e2e/f52b368f-d144-49b6-aed7-4c2624f1faca.md", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/84e6be187877360c5840ef63977ae3e9d7b193eb/e2e/811bc1ec-55f0-430f-803c-fcce64e8f840.md", "", "", "811bc1ec-55f0-430f-803c-fcce64e8f840.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/84e6be187877360c5840ef63977ae3e9d7b193eb/e2e/9b9a0341-eb6a-426a-864a-6d89105bbfa5.md", "", "", "9b9a0341-eb6a-426a-864a-6d89105bbfa5.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/da563445bb66940a7fb77f0f60f6d6413d89018f/e2e/21d5ef20-4a71-4957-bc19-8864fef6ab6f.md", "", "", "21d5ef20-4a71-4957-bc19-8864fef6ab6f.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/bf32a642862f2bbc4fa6fcf0526951adf1ca0d92/e2e/4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md", "", "", "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/e4317ec838fcb7d85ac106cbc52263e36e686afe/e2e/ac98edee-11e5-4200-8f9b-593df75e91d7.md", "", "", "ac98edee-11e5-4200-8f9b-593df75e91d7.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(5).Copy()
$ws2.Rows.Item(6).Insert()
$ws2.Rows.Item(5).Copy()
$ws2.Rows.Item(6).Insert()

$ws2.Range("A5").Value = "21d5ef20-4a71-4957-bc19-8864fef6ab6f.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "21d5ef20-4a71-4957-bc19-8864fef6ab6f.e5cf6278b4302378db34465a7e02f191588d849c.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-23 22:40:12"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("J5").Value = "Include"

$ws2.Range("A6").Value = "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md"
$ws2.Range("B6").Value = ".md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("D6").Value = "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.9b0eda3e4b066d8d5d1bc2b7f6c222923865dc25.zh-cn.xlf"
$ws2.Range("E6").Value = "2016-03-23 22:40:12"
$ws2.Range("H6").Value = "0001-01-01 00:00:00"
$ws2.Range("J6").Value = "Include"

# Row 7 keeps the duplicated ac98edee-... content untouched.

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a896c4a09a321d440863c6933282f9952176d25/e2e/f52b368f-d144-49b6-aed7-4c2624f1faca.md", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e1cbdda56e7e1707a7b942cb4b460d8b5e948690/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.zh-cn.xlf", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6b737258b042f65412ffc9653c1e0c256d1d9a15/e2e/f52b368f-d144-49b6-aed7-4c2624f1faca.md", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/215a79781b785b3f29b3e264fa47c7f16ab52970/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.zh-cn.xlf", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/84e6be187877360c5840ef63977ae3e9d7b193eb/e2e/811bc1ec-55f0-430f-803c-fcce64e8f840.md", "", "", "811bc1ec-55f0-430f-803c-fcce64e8f840.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5b4cd165dcfff63359d0734f5e1df93feded3da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/811bc1ec-55f0-430f-803c-fcce64e8f840.35db9e88cb7d9edfaedae33c7f1f64bfd40c182b.zh-cn.xlf", "", "", "811bc1ec-55f0-430f-803c-fcce64e8f840.35db9e88cb7d9edfaedae33c7f1f64bfd40c182b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/84e6be187877360c5840ef63977ae3e9d7b193eb/e2e/9b9a0341-eb6a-426a-864a-6d89105bbfa5.md", "", "", "9b9a0341-eb6a-426a-864a-6d89105bbfa5.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5b4cd165dcfff63359d0734f5e1df93feded3da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9b9a0341-eb6a-426a-864a-6d89105bbfa5.484cf389793b85af706d008733585fbea1189294.zh-cn.xlf", "", "", "9b9a0341-eb6a-426a-864a-6d89105bbfa5.484cf389793b85af706d008733585fbea1189294.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/da563445bb66940a7fb77f0f60f6d6413d89018f/e2e/21d5ef20-4a71-4957-bc19-8864fef6ab6f.md", "", "", "21d5ef20-4a71-4957-bc19-8864fef6ab6f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da563445bb66940a7fb77f0f60f6d6413d89018f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21d5ef20-4a71-4957-bc19-8864fef6ab6f.e5cf6278b4302378db34465a7e02f191588d849c.zh-cn.xlf", "", "", "21d5ef20-4a71-4957-bc19-8864fef6ab6f.e5cf6278b4302378db34465a7e02f191588d849c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/bf32a642862f2bbc4fa6fcf0526951adf1ca0d92/e2e/4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md", "", "", "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf32a642862f2bbc4fa6fcf0526951adf1ca0d92/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4999b8bf-edcd-4b0b-bbf6-e5582e2da185.9b0eda3e4b066d8d5d1bc2b7f6c222923865dc25.zh-cn.xlf", "", "", "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.9b0eda3e4b066d8d5d1bc2b7f6c222923865dc25.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/e4317ec838fcb7d85ac106cbc52263e36e686afe/e2e/ac98edee-11e5-4200-8f9b-593df75e91d7.md", "", "", "ac98edee-11e5-4200-8f9b-593df75e91d7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c8b98565e194f28ec4d09d15d0aa71076ff781d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ac98edee-11e5-4200-8f9b-593df75e91d7.9b57e2a4f996b2a240c1ecb33e07f1d7fb84ac7d.zh-cn.xlf", "", "", "ac98edee-11e5-4200-8f9b-593df75e91d7.9b57e2a4f996b2a240c1ecb33e07f1d7fb84ac7d.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(5).Copy()
$ws3.Rows.Item(6).Insert()
$ws3.Rows.Item(5).Copy()
$ws3.Rows.Item(6).Insert()

$ws3.Range("A5").Value = "21d5ef20-4a71-4957-bc19-8864fef6ab6f.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "21d5ef20-4a71-4957-bc19-8864fef6ab6f.e5cf6278b4302378db34465a7e02f191588d849c.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-23 22:40:20"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("J5").Value = "Include"

$ws3.Range("A6").Value = "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md"
$ws3.Range("B6").Value = ".md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("D6").Value = "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.9b0eda3e4b066d8d5d1bc2b7f6c222923865dc25.de-de.xlf"
$ws3.Range("E6").Value = "2016-03-23 22:40:20"
$ws3.Range("H6").Value = "0001-01-01 00:00:00"
$ws3.Range("J6").Value = "Include"

# Row 7 keeps the duplicated ac98edee-... content untouched.

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a896c4a09a321d440863c6933282f9952176d25/e2e/f52b368f-d144-49b6-aed7-4c2624f1faca.md", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1714f698b383539f0c0b62aa29e9a758aa910bf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.de-de.xlf", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7dab3f51353e50d40ad00660eda53c21eec84ce0/e2e/f52b368f-d144-49b6-aed7-4c2624f1faca.md", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/397ab83a3cc8e928042701b0d0740a57a22af6b8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.de-de.xlf", "", "", "f52b368f-d144-49b6-aed7-4c2624f1faca.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/84e6be187877360c5840ef63977ae3e9d7b193eb/e2e/811bc1ec-55f0-430f-803c-fcce64e8f840.md", "", "", "811bc1ec-55f0-430f-803c-fcce64e8f840.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7abd7031f46cf68af5037cb2b7bbc4e96b5e205b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/811bc1ec-55f0-430f-803c-fcce64e8f840.35db9e88cb7d9edfaedae33c7f1f64bfd40c182b.de-de.xlf", "", "", "811bc1ec-55f0-430f-803c-fcce64e8f840.35db9e88cb7d9edfaedae33c7f1f64bfd40c182b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/84e6be187877360c5840ef63977ae3e9d7b193eb/e2e/9b9a0341-eb6a-426a-864a-6d89105bbfa5.md", "", "", "9b9a0341-eb6a-426a-864a-6d89105bbfa5.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7abd7031f46cf68af5037cb2b7bbc4e96b5e205b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9b9a0341-eb6a-426a-864a-6d89105bbfa5.484cf389793b85af706d008733585fbea1189294.de-de.xlf", "", "", "9b9a0341-eb6a-426a-864a-6d89105bbfa5.484cf389793b85af706d008733585fbea1189294.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/da563445bb66940a7fb77f0f60f6d6413d89018f/e2e/21d5ef20-4a71-4957-bc19-8864fef6ab6f.md", "", "", "21d5ef20-4a71-4957-bc19-8864fef6ab6f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da563445bb66940a7fb77f0f60f6d6413d89018f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21d5ef20-4a71-4957-bc19-8864fef6ab6f.e5cf6278b4302378db34465a7e02f191588d849c.de-de.xlf", "", "", "21d5ef20-4a71-4957-bc19-8864fef6ab6f.e5cf6278b4302378db34465a7e02f191588d849c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/bf32a642862f2bbc4fa6fcf0526951adf1ca0d92/e2e/4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md", "", "", "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf32a642862f2bbc4fa6fcf0526951adf1ca0d92/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4999b8bf-edcd-4b0b-bbf6-e5582e2da185.9b0eda3e4b066d8d5d1bc2b7f6c222923865dc25.de-de.xlf", "", "", "4999b8bf-edcd-4b0b-bbf6-e5582e2da185.9b0eda3e4b066d8d5d1bc2b7f6c222923865dc25.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/e4317ec838fcb7d85ac106cbc52263e36e686afe/e2e/ac98edee-11e5-4200-8f9b-593df75e91d7.md", "", "", "ac98edee-11e5-4200-8f9b-593df75e91d7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b48c7d0e90f9e5b0dc7e8bb87e3de61989584bc4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ac98edee-11e5-4200-8f9b-593df75e91d7.9b57e2a4f996b2a240c1ecb33e07f1d7fb84ac7d.de-de.xlf", "", "", "ac98edee-11e5-4200-8f9b-593df75e91d7.9b57e2a4f996b2a240c1ecb33e07f1d7fb84ac7d.de-de.xlf") | Out-Null

Write-Host "Report generated for handoff: added 21d5ef20-... and 4999b8bf-... rows."
